# "communiaction db -> crawler": the requests table's `date` column values
# are bumped forward (communication DB dates -> crawler run dates), and the
# active sheet/selection moves from "cities" to "requests" (cell E4).

$wb = $excel.ActiveWorkbook

# --- Update the "requests" sheet's date column (E2:E4) ---
$requests = $wb.Worksheets.Item("requests")

$requests.Range("E2").Value = 43832
$requests.Range("E3").Value = 43832
$requests.Range("E4").Value = 43833

# --- Move the active sheet/selection to "requests"!E4 ---
# (this also clears tabSelected/selection bookkeeping on whichever sheet
# used to be active, e.g. "cities")
$requests.Activate() | Out-Null
$requests.Range("E4").Select() | Out-Null
